$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = '22.446.67'
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = '  +0.13%  '
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = '1.573.38'
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = '  +0.07%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '287.80'
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = '  +0.87%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.3704'
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = '  +1.84%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '47.27'
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = '  -2.69%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.3318'
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = '  -0.62%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '1.156'
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = '  +2.46%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.07504'
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = '  +0.06%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '20.75'
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = '  -0.38%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '5.933'
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = '  +0.00%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '6.919'
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = '  +0.38%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '1.564.19'
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = '  -0.01%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.00001116'
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = '  +0.94%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '88.36'
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = '  +0.18%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '0.06726'
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = '  +0.60%  '
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = 'Dai'
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '1.001'
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = '  +0.09%  '
$ws.Range("B21").NumberFormat = "@"
$ws.Range("B21").Value = 'Uniswap'
$ws.Range("C21").NumberFormat = "@"
$ws.Range("C21").Value = 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni'
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '6.397'
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = '  +0.51%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '16.50'
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = '  +2.35%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '11.98'
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = '  -0.26%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '22.437.86'
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = '  +0.12%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '2.379'
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = '  -1.25%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '2.627'
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = '  +2.66%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '150.88'
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = '  +0.88%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '19.61'
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = '  +0.96%  '
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = '  -0.69%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '124.97'
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = '  +1.49%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '1.741.73'
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = '  +0.15%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '1.094'
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = '  +2.89%  '
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = '  -1.08%  '
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = '  -0.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '9.908'
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = '  +0.77%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.08330'
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = '  +0.99%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02443'
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = '  +1.85%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.304'
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = '  -0.15%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.06379'
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = '  -0.13%  '
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = '  +0.32%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '5.329'
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = '  -0.18%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.6235'
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = '  +2.44%  '
$ws.Range("B44").NumberFormat = "@"
$ws.Range("B44").Value = 'EnergySwap'
$ws.Range("C44").NumberFormat = "@"
$ws.Range("C44").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '14.05'
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = '  +1.60%  '
$ws.Range("B45").NumberFormat = "@"
$ws.Range("B45").Value = 'Decentraland'
$ws.Range("C45").NumberFormat = "@"
$ws.Range("C45").Value = 'https://coinranking.com/coin/tEf7-dnwV3BXS+decentraland-mana'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.6047'
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = '  +4.91%  '
$ws.Range("B46").NumberFormat = "@"
$ws.Range("B46").Value = 'PancakeSwap'
$ws.Range("C46").NumberFormat = "@"
$ws.Range("C46").Value = 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '3.773'
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = '  +0.44%  '
$ws.Range("B47").NumberFormat = "@"
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").NumberFormat = "@"
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '2.041'
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = '  +1.40%  '
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = 'Quant'
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '124.76'
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = '  +0.02%  '
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = 'EOS'
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = 'https://coinranking.com/coin/iAzbfXiBBKkR6+eos-eos'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '1.206'
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = '  -0.82%  '
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = 'Cronos'
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.07195'
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = '  -0.27%  '
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = 'Aave'
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = 'https://coinranking.com/coin/ixgUfzmLR+aave-aave'
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '77.15'
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = '  +2.37%  '
